$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated parameter-estimation results (row 1, columns A:Q).
$ws.Range("A1").Value = -0.00483153360588449
$ws.Range("B1").Value = 0.00439763927767704
$ws.Range("C1").Value = 0.000468599539723256
$ws.Range("D1").Value = 0.000501779565624479
$ws.Range("E1").Value = 2.49495984176715
$ws.Range("F1").Value = 0.00145780303687494
$ws.Range("G1").Value = 0.148273757923292
$ws.Range("H1").Value = 0.0283184279154154
$ws.Range("I1").Value = 0.318058863281476
$ws.Range("J1").Value = 2.37907604454986
$ws.Range("K1").Value = 0.0481521631794444
$ws.Range("L1").Value = 0.28312955275417
$ws.Range("M1").Value = 0.0949889760089385
$ws.Range("N1").Value = 0.00646523353368974
$ws.Range("O1").Value = 0.000842791378337106
$ws.Range("P1").Value = 0.0584393786336576
$ws.Range("Q1").Value = -0.0312982736170995

# C1, D1 and O1 previously carried a custom scientific-notation number format;
# that direct formatting is removed so the cells fall back to the default style.
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").Style = "Normal"
$ws.Range("O1").Style = "Normal"

# The saved selection now spans the whole data row instead of just A2.
$ws.Range("A1:Q1").Select()
